$d = $word.ActiveDocument

$d.Content.Find.Execute("Я, я лучший по", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Я, Соколова Наталья Михайловна", 2)

$d.Content.Find.Execute(" разрешаю обучающемуся коркина соня", $true, $false, $false, $false, $false,
                         $true, 1, $false, " разрешаю обучающемуся Пупкин Василий", 2)

$d.Content.Find.Execute(" 4 урока по причине: у меня олимпиада", $true, $false, $false, $false, $false,
                         $true, 1, $false, " 5 урока по причине: хочу есть", 2)

$d.Content.Find.Execute("Классный руководитель    2022-11-14 (дата)                                      ___________ (подпись)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Классный руководитель    2022-11-17 (дата)                                      ___________ (подпись)", 2)

$d.Content.Find.Execute("Представитель администрации 2022-11-14 (дата)                           ___________ (подпись)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Представитель администрации 2022-11-17 (дата)                           ___________ (подпись)", 2)
